$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 'Hapoel TelAviv'

# Row 7
$ws.Range("F7").Value = 'Hapoel Haifa'

# Row 10
$ws.Range("B10").Value = 6799825
$ws.Range("F10").Value = 'Maccabi Petach Tikva'
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 1
$ws.Range("K10").Value = 'D'
$ws.Range("L10").Value = 1.8
$ws.Range("M10").Value = 3.25
$ws.Range("N10").Value = 4
$ws.Range("O10").Value = 1.95
$ws.Range("P10").Value = 3.2
$ws.Range("Q10").Value = 3.4
$ws.Range("R10").Value = -0.5
$ws.Range("S10").Value = 2.1
$ws.Range("T10").Value = 1.775
$ws.Range("V10").Value = 2.025
$ws.Range("W10").Value = 1.825
$ws.Range("X10").Value = -1
$ws.Range("Y10").Value = 2.2
$ws.Range("AA10").Value = -1
$ws.Range("AB10").Value = 0.7749999999999999
$ws.Range("AC10").Value = 1.025
$ws.Range("AD10").Value = -1

# Row 11
$ws.Range("B11").Value = 6799829
$ws.Range("F11").Value = 'Maccabi Netanya'
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 'H'
$ws.Range("L11").Value = 2.4
$ws.Range("M11").Value = 3.4
$ws.Range("N11").Value = 2.6
$ws.Range("O11").Value = 2.625
$ws.Range("P11").Value = 3.4
$ws.Range("Q11").Value = 2.375
$ws.Range("R11").Value = 0
$ws.Range("S11").Value = 2
$ws.Range("T11").Value = 1.85
$ws.Range("V11").Value = 2
$ws.Range("W11").Value = 1.85
$ws.Range("X11").Value = 1.625
$ws.Range("Y11").Value = -1
$ws.Range("AA11").Value = 1
$ws.Range("AB11").Value = -1
$ws.Range("AC11").Value = -1
$ws.Range("AD11").Value = 0.8500000000000001

# Row 16
$ws.Range("E16").Value = 'Hapoel Haifa'

# Row 21
$ws.Range("F21").Value = 'Hapoel TelAviv'

# Row 23
$ws.Range("F23").Value = 'Hapoel Haifa'

# Row 24
$ws.Range("B24").Value = 6799841
$ws.Range("E24").Value = 'MS Ashdod'
$ws.Range("F24").Value = 'Maccabi Petach Tikva'
$ws.Range("L24").Value = 2.25
$ws.Range("M24").Value = 3.25
$ws.Range("N24").Value = 2.75
$ws.Range("O24").Value = 2.3
$ws.Range("P24").Value = 3.3
$ws.Range("Q24").Value = 2.7
$ws.Range("S24").Value = 1.75
$ws.Range("T24").Value = 2.05
$ws.Range("V24").Value = 2
$ws.Range("W24").Value = 1.85
$ws.Range("Y24").Value = 2.3
$ws.Range("AD24").Value = 0.8500000000000001

# Row 25
$ws.Range("B25").Value = 6799846
$ws.Range("E25").Value = 'Hapoel Jerusalem FC'
$ws.Range("F25").Value = 'Maccabi Netanya'
$ws.Range("L25").Value = 2.8
$ws.Range("M25").Value = 3.3
$ws.Range("N25").Value = 2.3
$ws.Range("O25").Value = 2.5
$ws.Range("P25").Value = 3.2
$ws.Range("Q25").Value = 2.6
$ws.Range("S25").Value = 1.9
$ws.Range("T25").Value = 1.95
$ws.Range("V25").Value = 2.05
$ws.Range("W25").Value = 1.8
$ws.Range("Y25").Value = 2.2
$ws.Range("AD25").Value = 0.8

# Row 27
$ws.Range("E27").Value = 'Hapoel TelAviv'

# Row 31
$ws.Range("E31").Value = 'Hapoel Haifa'

# Row 35
$ws.Range("F35").Value = 'Hapoel TelAviv'

# Row 37
$ws.Range("E37").Value = 'Hapoel TelAviv'
$ws.Range("F37").Value = 'Hapoel Haifa'

# Row 46
$ws.Range("E46").Value = 'Hapoel Haifa'

# Row 47
$ws.Range("F47").Value = 'Hapoel TelAviv'

# Row 52
$ws.Range("E52").Value = 'Hapoel TelAviv'

# Row 53
$ws.Range("F53").Value = 'Hapoel Haifa'

# Row 57
$ws.Range("F57").Value = 'Hapoel TelAviv'

# Row 63
$ws.Range("E63").Value = 'Hapoel Haifa'

# Row 66
$ws.Range("E66").Value = 'Hapoel TelAviv'

# Row 67
$ws.Range("F67").Value = 'Hapoel Haifa'

# Row 71
$ws.Range("E71").Value = 'Hapoel TelAviv'

# Row 74
$ws.Range("E74").Value = 'Hapoel Haifa'

# Row 79
$ws.Range("F79").Value = 'Hapoel TelAviv'

# Row 80
$ws.Range("F80").Value = 'Hapoel Haifa'

# Row 90
$ws.Range("E90").Value = 'Hapoel Haifa'

# Row 91
$ws.Range("E91").Value = 'Hapoel TelAviv'

# Row 95
$ws.Range("E95").Value = 'Hapoel TelAviv'

# Row 98
$ws.Range("E98").Value = 'Hapoel Haifa'

# Row 101
$ws.Range("F101").Value = 'Hapoel TelAviv'

# Row 104
$ws.Range("F104").Value = 'Hapoel Haifa'

# Row 110
$ws.Range("F110").Value = 'Hapoel Haifa'

# Row 111
$ws.Range("E111").Value = 'Hapoel TelAviv'

# Row 114
$ws.Range("E114").Value = 'Hapoel Haifa'

# Row 120
$ws.Range("F120").Value = 'Hapoel TelAviv'

# Row 122
$ws.Range("F122").Value = 'Hapoel Haifa'

# Row 127
$ws.Range("E127").Value = 'Hapoel TelAviv'

# Row 132
$ws.Range("E132").Value = 'Hapoel Haifa'
$ws.Range("F132").Value = 'Hapoel TelAviv'

# Row 137
$ws.Range("F137").Value = 'Hapoel Haifa'

# Row 138
$ws.Range("E138").Value = 'Hapoel TelAviv'

# Row 143
$ws.Range("F143").Value = 'Hapoel TelAviv'

# Row 148
$ws.Range("E148").Value = 'Hapoel Haifa'

# Row 149
$ws.Range("E149").Value = 'Hapoel TelAviv'

# Row 154
$ws.Range("F154").Value = 'Hapoel Haifa'

# Row 156
$ws.Range("F156").Value = 'Hapoel TelAviv'

# Row 157
$ws.Range("E157").Value = 'Hapoel Haifa'

# Row 163
$ws.Range("F163").Value = 'Hapoel Haifa'

# Row 168
$ws.Range("F168").Value = 'Hapoel TelAviv'

# Row 173
$ws.Range("E173").Value = 'Hapoel TelAviv'

# Row 176
$ws.Range("E176").Value = 'Hapoel Haifa'

# Row 182
$ws.Range("F182").Value = 'Hapoel Haifa'

# Row 183
$ws.Range("F183").Value = 'Hapoel TelAviv'

# Row 184
$ws.Range("F184").Value = 'Hapoel Haifa'

# Row 190
$ws.Range("F190").Value = 'Hapoel TelAviv'

# Row 191
$ws.Range("F191").Value = 'Hapoel Haifa'

# Row 196
$ws.Range("E196").Value = 'Hapoel TelAviv'

# Row 198
$ws.Range("E198").Value = 'Hapoel Haifa'

# Row 201
$ws.Range("F201").Value = 'Hapoel Haifa'

# Row 202
$ws.Range("B202").Value = 8016168
$ws.Range("E202").Value = 'Maccabi Netanya'
$ws.Range("F202").Value = 'Hapoel Jerusalem FC'
$ws.Range("G202").Value = 1
$ws.Range("H202").Value = 0
$ws.Range("I202").Value = 1
$ws.Range("K202").Value = 'H'
$ws.Range("L202").Value = 2
$ws.Range("M202").Value = 3.3
$ws.Range("N202").Value = 3.6
$ws.Range("O202").Value = 2.1
$ws.Range("P202").Value = 3.3
$ws.Range("Q202").Value = 3.4
$ws.Range("R202").Value = -0.25
$ws.Range("S202").Value = 1.85
$ws.Range("T202").Value = 2
$ws.Range("U202").Value = 2.25
$ws.Range("V202").Value = 1.95
$ws.Range("W202").Value = 1.9
$ws.Range("X202").Value = 1.1
$ws.Range("Z202").Value = -1
$ws.Range("AA202").Value = 0.8500000000000001
$ws.Range("AB202").Value = -1
$ws.Range("AC202").Value = -1
$ws.Range("AD202").Value = 0.8999999999999999

# Row 203
$ws.Range("B203").Value = 8015667
$ws.Range("E203").Value = 'Hapoel Bnei Sakhnin'
$ws.Range("F203").Value = 'Maccabi Bnei Raina'
$ws.Range("G203").Value = 0
$ws.Range("H203").Value = 2
$ws.Range("I203").Value = 0
$ws.Range("K203").Value = 'A'
$ws.Range("L203").Value = 2.6
$ws.Range("M203").Value = 3.1
$ws.Range("N203").Value = 2.75
$ws.Range("O203").Value = 2.45
$ws.Range("P203").Value = 3
$ws.Range("Q203").Value = 2.9
$ws.Range("R203").Value = 0
$ws.Range("S203").Value = 1.775
$ws.Range("T203").Value = 2.1
$ws.Range("U203").Value = 2
$ws.Range("V203").Value = 1.825
$ws.Range("W203").Value = 2.025
$ws.Range("X203").Value = -1
$ws.Range("Z203").Value = 1.9
$ws.Range("AA203").Value = -1
$ws.Range("AB203").Value = 1.1
$ws.Range("AC203").Value = 0
$ws.Range("AD203").Value = 0

# Row 204
$ws.Range("B204").Value = 8016167
$ws.Range("E204").Value = 'Maccabi Petach Tikva'
$ws.Range("F204").Value = 'Hapoel TelAviv'
$ws.Range("H204").Value = 2
$ws.Range("L204").Value = 2.75
$ws.Range("M204").Value = 2.8
$ws.Range("N204").Value = 2.875
$ws.Range("O204").Value = 3.25
$ws.Range("P204").Value = 2.875
$ws.Range("Q204").Value = 2.45
$ws.Range("R204").Value = 0.25
$ws.Range("S204").Value = 1.825
$ws.Range("T204").Value = 2.025
$ws.Range("V204").Value = 1.975
$ws.Range("W204").Value = 1.875
$ws.Range("X204").Value = 2.25
$ws.Range("AA204").Value = 0.825
$ws.Range("AC204").Value = 0.9750000000000001

# Row 205
$ws.Range("B205").Value = 8015668
$ws.Range("E205").Value = 'Maccabi Haifa'
$ws.Range("F205").Value = 'Hapoel Beer Sheva'
$ws.Range("H205").Value = 1
$ws.Range("L205").Value = 1.833
$ws.Range("M205").Value = 3.6
$ws.Range("N205").Value = 4
$ws.Range("O205").Value = 1.833
$ws.Range("P205").Value = 3.5
$ws.Range("Q205").Value = 4.2
$ws.Range("R205").Value = -0.5
$ws.Range("S205").Value = 1.875
$ws.Range("T205").Value = 1.975
$ws.Range("V205").Value = 1.85
$ws.Range("W205").Value = 2
$ws.Range("X205").Value = 0.833
$ws.Range("AA205").Value = 0.875
$ws.Range("AC205").Value = 0.8500000000000001

# Row 211
$ws.Range("E211").Value = 'Hapoel TelAviv'

# Row 214
$ws.Range("E214").Value = 'Hapoel Haifa'

# Row 215
$ws.Range("F215").Value = 'Hapoel TelAviv'

# Row 216
$ws.Range("E216").Value = 'Hapoel Haifa'

# Row 222
$ws.Range("E222").Value = 'Hapoel Haifa'

# Row 227
$ws.Range("E227").Value = 'Hapoel TelAviv'

# Row 230
$ws.Range("F230").Value = 'Hapoel TelAviv'

# Row 232
$ws.Range("F232").Value = 'Hapoel Haifa'

# Row 236
$ws.Range("E236").Value = 'Hapoel Haifa'

# Row 239
$ws.Range("F239").Value = 'Hapoel Haifa'
